$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Row 4: Cole H. testimonial ---
$tbl.ListRows.Add() | Out-Null
$ws.Range("A3:K3").Copy()
$ws.Range("A4:K4").PasteSpecial(-4122)

$ws.Cells.Item(4,1).Value = 45984.67118607639
$ws.Cells.Item(4,2).Value = "chunzeker24@gmail.com"
$ws.Cells.Item(4,3).Value = "The biggest challenge I was facing at the time was lack of upward mobility, and preparedness if the right opportunity came up."
$ws.Cells.Item(4,4).Value = "I felt stuck and was unsure how to get to the next level of my career. At times it would lead to burnout in my current role because it felt like I would never get another opportunity."
$ws.Cells.Item(4,5).Value = "The most helpful advice I received was in regards to positioning myself for the next step. I was told to go above and beyond in my current role. That didn’t just mean in the work I do, but also the connections I make. Additionally, I was mentored on proper interviewing skills for when I do get an opportunity."
$ws.Cells.Item(4,6).Value = "I believe you are very honest and have went through similar experiences as I have. You were specific in your advice and even gave real life suggestions pertaining to my job that would help with career advancement. You know it’s not a one size fits all for advancement and your ability to recognize that and give me the ingredients specific to my aspirations l was very helpful. "
$ws.Cells.Item(4,7).Value = "I received an interview for a promotion. While I didn’t get the job, I got great feedback from the interviewer about my interviewing skills and the only reason I didn’t get the job was there was someone already in line for the position. I have also gotten more recognition on my team as a result of implementing the advice I was given."
$ws.Cells.Item(4,8).Value = "It’s hard to put a number on the help received but I think it kickstarted my advancement. I believe I will get a promotion a year or 2 earlier than I would have without the help and advice. I now feel positioned well for the future. "
$ws.Cells.Item(4,9).Value = "I would say there are many reasons. You are always willing to lend advice, many times based on personal experiences. You also have a history of success and advancement which is always a positive when working with another individual. Most importantly, you are honest and adaptable which helps with framing advice for an individual. "
$ws.Cells.Item(4,10).Value = "Accelerate: getting a promotion or a leadership position"
$ws.Cells.Item(4,11).Value = "Yes"

# --- Row 5: Ebuka O. testimonial ---
$tbl.ListRows.Add() | Out-Null
$ws.Range("A3:K3").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)

$ws.Cells.Item(5,1).Value = 45984.78307002315
$ws.Cells.Item(5,2).Value = "chiebukaonyejesi@gmail.com"
$ws.Cells.Item(5,3).Value = "I was starting out as an intern at John Deere and needed to know how to maximize my internship to get a return offer and overall make an impact at the company."
$ws.Cells.Item(5,4).Value = "I felt kinda lost and it was something I wanted to overcome at the time. I felt like getting a coach would make it easier to cross that barrier and achieve my plans for the summer faster. "
$ws.Cells.Item(5,5).Value = "You told me not to put myself in a box and just limit my experience/career to John Deere. That's something I was shocked to hear because most people will try to convince you to work for their company but you said there's a lot of opportunities everywhere and I should really try to find what really resonates with me."
$ws.Cells.Item(5,6).Value = "I felt like you were more raw and didn’t just say generic stuff. You said things and gave feedback that was actually applicable and not just fluff. I like how you also didn’t just jump into giving advice you kind of dug a bit deeper before your responses."
$ws.Cells.Item(5,7).Value = "I believe I connected with more people at John Deere and also made the most out of my intern summer project because I was in direct contact with my manger. I also started thinking more about my career out of John Deere and continue working fully on the app I was developing because you highlighted everyone's path is different."
$ws.Cells.Item(5,8).Value = "I would say there has been infact an impact in my career trajectory. I have had some changes in mindset. I think more deeply and broadly about things related to my career now and try not to limit myself or get stuck up on a specific role/company."
$ws.Cells.Item(5,9).Value = "As I said earlier, he's very raw and wouldn't just tell you things you would like to hear. I believe he would give you realistic and actionable advice and will also adapt it to your current situation."
$ws.Cells.Item(5,10).Value = "Launch: getting the first job or an entry level job"
$ws.Cells.Item(5,11).Value = "Yes"

$ws.Rows.Item(4).RowHeight = 22.5
$ws.Rows.Item(5).RowHeight = 22.5

# Keep the hidden _FilterDatabase defined name in sync with the expanded table range
$filterDbName = $wb.Names.Item(1)
$filterDbName.RefersTo = "='Form Responses 1'!`$A`$1:`$K`$5"
